$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet1" to "Estimates".
# (This also updates the _xlnm._FilterDatabase defined names that reference the sheet.)
$ws.Name = "Estimates"

# Row 65, column A currently reuses the "Min (P=95%)" label (shared with row 64).
# Give it its own distinct label "Max (P=95%)" since it represents the upper bound.
$ws.Range("A65").Value = "Max (P=95%)"
